$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy header style/format from an existing header cell (A1) onto the new headers
$ws.Range("A1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)

# Boolean "Outliers_MAD" columns
$ws.Range("F2").Value = $true
$ws.Range("G2").Value = $false
$ws.Range("H2").Value = $false

$ws.Range("F3").Value = $false
$ws.Range("G3").Value = $false
$ws.Range("H3").Value = $false

$ws.Range("F4").Value = $false
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = $false

$ws.Range("F5").Value = $false
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = $false

$ws.Range("F6").Value = $true
$ws.Range("G6").Value = $false
$ws.Range("H6").Value = $false
